# Coupling Parameters sheet is the active sheet (rId2 / sheet2.xml)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Add a new "InvestmentIteration" parameter row below the existing ones,
# used to persist project results across coupling/investment iterations.
$ws.Range("A8").Value = "InvestmentIteration"
$ws.Range("B8").Value = 0

# Mirror the cursor position left behind in the saved workbook.
$ws.Range("D9").Select()
